$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column V with header date "17-03-2020", copying the format (style) of the
# preceding header cell (U1) so it matches the other bold/bordered header cells.
$ws.Range("U1").Copy()
$ws.Range("V1").PasteSpecial(-4122)
$ws.Range("V1").Value = "17-03-2020"

# Fill in the new confirmed-case counts per province for 17-03-2020
$ws.Range("V2").Value = 17
$ws.Range("V3").Value = 24
$ws.Range("V4").Value = 14
$ws.Range("V5").Value = 173
$ws.Range("V6").Value = 10
$ws.Range("V7").Value = 197
$ws.Range("V8").Value = 634
$ws.Range("V9").Value = 151
$ws.Range("V10").Value = 45
$ws.Range("V11").Value = 173
$ws.Range("V12").Value = 20
$ws.Range("V13").Value = 175

$ws.Range("A1").Select()
